$d = $word.ActiveDocument

# Target the first paragraph (the "**ID__...__ID**" marker paragraph).
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right), 5pt space-from-text on each
# side, matching <w:pBdr><w:top w:space="5"/><w:left .../><w:bottom .../>
# <w:right .../></w:pBdr> with no explicit line style/size/color.
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Change the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Remove the trailing run that contains only a literal space
# (<w:r><w:rPr>...</w:rPr><w:t xml:space="preserve"> </w:t></w:r>),
# leaving the paragraph mark intact. The paragraph mark is the last
# character of the paragraph's range, and the space is the character
# immediately before it.
$spaceRange = $d.Range($p1.Range.End - 2, $p1.Range.End - 1)
$spaceRange.Text = ""

# Update the marker text in the (now sole) run of the paragraph.
$p1.Range.Find.Execute("**ID__AFFARS_5341_topic_5__ID**", $true, $false, $false, $false, $false,
                        $true, 1, $false, "**ID__AFFARS_5341_202__ID**", 2)
